$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-11 (columns A-E)
$data = @(
    @(0, 0.4019, 316, 1790, 0.001315),
    @(1, 0.3409, 173, 949,  0.000736),
    @(2, 0.316,  129, 697,  0.000557),
    @(3, 0.1698, 153, 866,  0.000615),
    @(4, 0.1285, 205, 1167, 0.00081),
    @(5, 0.1715, 213, 1206, 0.000837),
    @(6, 0.1154, 241, 1364, 0.000945),
    @(7, 0.06819,257, 1472, 0.001038),
    @(8, 0.05369, 91, 504,  0.000393),
    @(9, 0.04,   152, 839,  0.000596)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
}

# Remove the now-obsolete last row (row 12) entirely, shrinking the used range
$ws.Rows.Item(12).Delete()
